$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Vecka: 15 -> 16
$ws.Range("B3").Value = 16

# Shift the date column forward by three weeks (2015-03-23.. -> 2015-04-13..)
$ws.Range("A9").Value  = [datetime]"2015-04-13"
$ws.Range("A10").Value = [datetime]"2015-04-14"
$ws.Range("A11").Value = [datetime]"2015-04-15"
$ws.Range("A12").Value = [datetime]"2015-04-16"
$ws.Range("A13").Value = [datetime]"2015-04-17"
$ws.Range("A14").Value = [datetime]"2015-04-18"
$ws.Range("A15").Value = [datetime]"2015-04-19"

# Activity notes for the week
$ws.Range("B10").Value = "Puzzeldesign och nivådesign"
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = "nivådesign"
$ws.Range("B13").Value = "nivådesign"
$ws.Range("B14").Value = "nivådesign"
$ws.Range("B15").Value = "nivådesign"

# Hours logged per day
$ws.Range("D12").Value = 2
$ws.Range("D13").Value = 4
$ws.Range("D14").Value = 3
$ws.Range("D15").Value = 2

# Resize column A to fit the new date values
$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("A15").Select() | Out-Null
